# Fix crear caja estimacion en presupuesto
# Updates values in the estimation/cash-box sheet to reflect corrected figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 4476000

$ws.Range("B7").Value = 4488500

$ws.Range("B9").Value = 8952
$ws.Range("C9").Value = 8952

$ws.Range("B12").Value = 12252
$ws.Range("C12").Value = 12252

$ws.Range("B13").Value = 4476248
$ws.Range("C13").Value = -12252
